# Apply updated odds values to Jogos_da_Semana_FlashScore_2025-04-02 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.36
$ws.Range("H2").Value = 5.25
$ws.Range("I2").Value = 7.5
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 1.2
$ws.Range("M2").Value = 4.5
$ws.Range("N2").Value = 1.67
$ws.Range("O2").Value = 2.2
$ws.Range("P2").Value = 1.3
$ws.Range("Q2").Value = 3.4
$ws.Range("R2").Value = 1.95
$ws.Range("S2").Value = 1.8
$ws.Range("T2").Value = 8
$ws.Range("U2").Value = 7
$ws.Range("W2").Value = 9.5
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 13
$ws.Range("AB2").Value = 23
$ws.Range("AC2").Value = 67
$ws.Range("AD2").Value = 451
$ws.Range("AG2").Value = 26
$ws.Range("AH2").Value = 101
$ws.Range("AI2").Value = 67
# Row 3
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 3.2
$ws.Range("L3").Value = 1.22
$ws.Range("M3").Value = 4.33
$ws.Range("N3").Value = 1.73
$ws.Range("O3").Value = 2.1
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 2.2
$ws.Range("Y3").Value = 26
$ws.Range("Z3").Value = 13
$ws.Range("AA3").Value = 7.5
$ws.Range("AB3").Value = 15
# Row 4
$ws.Range("G4").Value = 1.17
$ws.Range("H4").Value = 7.5
$ws.Range("I4").Value = 11
$ws.Range("N4").Value = 1.37
$ws.Range("U4").Value = 8
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 11
$ws.Range("Z4").Value = 21
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 67
$ws.Range("AD4").Value = 401
$ws.Range("AE4").Value = 34
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 34
$ws.Range("AH4").Value = 201
$ws.Range("AJ4").Value = 67
# Row 5
$ws.Range("G5").Value = 1.69
$ws.Range("N5").Value = 1.41
$ws.Range("O5").Value = 2.7
# Row 6
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 1.54
$ws.Range("L6").Value = 1.25
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 1.76
$ws.Range("O6").Value = 1.96
$ws.Range("T6").Value = 17
$ws.Range("Z6").Value = 12
$ws.Range("AB6").Value = 19
$ws.Range("AE6").Value = 8
# Row 7
$ws.Range("G7").Value = 1.45
$ws.Range("N7").Value = 1.93
$ws.Range("O7").Value = 1.97
$ws.Range("P7").Value = 1.36
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 1.75
$ws.Range("U7").Value = 7.5
$ws.Range("AF7").Value = 34
# Row 8
$ws.Range("G8").Value = 1.41
$ws.Range("H8").Value = 4.5
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 11
$ws.Range("L8").Value = 1.29
$ws.Range("M8").Value = 3.5
$ws.Range("N8").Value = 1.9
$ws.Range("O8").Value = 1.95
$ws.Range("W8").Value = 9.5
$ws.Range("X8").Value = 13
$ws.Range("Z8").Value = 11
$ws.Range("AA8").Value = 8.5
# Row 11
$ws.Range("G11").Value = 3.6
$ws.Range("I11").Value = 2.05
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1.75
$ws.Range("T11").Value = 9
$ws.Range("V11").Value = 13
$ws.Range("X11").Value = 34
$ws.Range("AC11").Value = 67
$ws.Range("AF11").Value = 9
$ws.Range("AH11").Value = 17
# Row 12
$ws.Range("G12").Value = 3.1
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 2.25
$ws.Range("T12").Value = 12
$ws.Range("U12").Value = 17
$ws.Range("V12").Value = 11
$ws.Range("W12").Value = 34
$ws.Range("X12").Value = 21
$ws.Range("AA12").Value = 7
$ws.Range("AD12").Value = 126
$ws.Range("AE12").Value = 10
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 9
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 17
$ws.Range("AJ12").Value = 23
# Row 13
$ws.Range("G13").Value = 2.88
$ws.Range("H13").Value = 3.2
$ws.Range("I13").Value = 2.5
$ws.Range("J13").Value = 1.05
$ws.Range("K13").Value = 11
$ws.Range("L13").Value = 1.25
$ws.Range("M13").Value = 3.75
$ws.Range("N13").Value = 1.88
$ws.Range("O13").Value = 1.98
$ws.Range("R13").Value = 1.67
$ws.Range("S13").Value = 2.1
$ws.Range("T13").Value = 10
$ws.Range("U13").Value = 15
$ws.Range("V13").Value = 11
$ws.Range("W13").Value = 29
$ws.Range("X13").Value = 21
$ws.Range("Z13").Value = 11
$ws.Range("AA13").Value = 6
$ws.Range("AF13").Value = 13
$ws.Range("AG13").Value = 10
$ws.Range("AH13").Value = 23
$ws.Range("AI13").Value = 19
$ws.Range("AJ13").Value = 26
# Row 14
$ws.Range("G14").Value = 1.29
$ws.Range("H14").Value = 7
$ws.Range("I14").Value = 7.5
$ws.Range("J14").Value = 1.01
$ws.Range("K14").Value = 34
$ws.Range("L14").Value = 1.04
$ws.Range("M14").Value = 13
$ws.Range("N14").Value = 1.17
$ws.Range("O14").Value = 5
$ws.Range("P14").Value = 1.13
$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = 1.36
$ws.Range("S14").Value = 3
$ws.Range("T14").Value = 21
$ws.Range("U14").Value = 13
$ws.Range("V14").Value = 11
$ws.Range("W14").Value = 13
$ws.Range("X14").Value = 10
$ws.Range("AA14").Value = 17
$ws.Range("AB14").Value = 17
$ws.Range("AC14").Value = 34
$ws.Range("AE14").Value = 41
$ws.Range("AG14").Value = 26
$ws.Range("AH14").Value = 81
# Row 15
$ws.Range("G15").Value = 2.15
$ws.Range("R15").Value = 2.2
$ws.Range("S15").Value = 1.62
$ws.Range("U15").Value = 9
$ws.Range("W15").Value = 21
$ws.Range("AB15").Value = 19
# Row 26
$ws.Range("G26").Value = 1.17
$ws.Range("H26").Value = 8
$ws.Range("J26").Value = 1.01
$ws.Range("K26").Value = 26
$ws.Range("L26").Value = 1.1
$ws.Range("M26").Value = 7
$ws.Range("N26").Value = 1.33
$ws.Range("O26").Value = 3.4
$ws.Range("P26").Value = 1.2
$ws.Range("Q26").Value = 4.33
$ws.Range("T26").Value = 10
$ws.Range("U26").Value = 7.5
$ws.Range("W26").Value = 7.5
$ws.Range("Y26").Value = 26
$ws.Range("Z26").Value = 23
$ws.Range("AB26").Value = 29
$ws.Range("AD26").Value = 351
# Row 27
$ws.Range("J27").Value = 1.07
$ws.Range("K27").Value = 7.5
# Row 28
$ws.Range("L28").Value = 1.22
$ws.Range("M28").Value = 4
$ws.Range("N28").Value = 1.73
$ws.Range("O28").Value = 2.08
# Row 29
$ws.Range("G29").Value = 1.6
$ws.Range("I29").Value = 5
$ws.Range("L29").Value = 1.22
$ws.Range("M29").Value = 4
$ws.Range("N29").Value = 1.75
$ws.Range("O29").Value = 2.05
$ws.Range("P29").Value = 1.33
$ws.Range("Q29").Value = 3.25
$ws.Range("R29").Value = 1.8
$ws.Range("S29").Value = 1.91
$ws.Range("T29").Value = 7.5
$ws.Range("U29").Value = 8
$ws.Range("W29").Value = 12
$ws.Range("Y29").Value = 23
$ws.Range("Z29").Value = 13
$ws.Range("AB29").Value = 17
$ws.Range("AD29").Value = 251
$ws.Range("AE29").Value = 15
$ws.Range("AG29").Value = 15
# Row 30
$ws.Range("I30").Value = 3.25
$ws.Range("J30").Value = 1.06
$ws.Range("K30").Value = 10
$ws.Range("N30").Value = 1.98
$ws.Range("O30").Value = 1.88
$ws.Range("U30").Value = 10
$ws.Range("W30").Value = 19
